# Add a new "Reg Proc" column (column T) to the "Table2" table on the
# "Details" sheet, populate its header/data cells, update one existing
# cell (S8) and move the current selection, matching the target revision.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Details")
$ws.Activate()

# --- Extend the table so column T becomes part of Table2 -------------------
$lo = $ws.ListObjects.Item("Table2")
$lo.Resize($ws.Range("A2:T16"))

# --- Header (row 2) ----------------------------------------------------
# Copy S2's header formatting (centered, bordered header style) onto T2,
# then overwrite with the new header text.
$ws.Range("S2").Copy($ws.Range("T2"))
$ws.Range("T2").Value = "Reg Proc"

# --- New column data cells (rows 4-10) ---------------------------------
$ws.Range("T4").Value = "When UIN IS needed to be generated`n1.the Acknowledgment from Print queue- what needs to be done`nTime period `n2. If there is a print failure- no need to handle from MOSIP`nUser Story ?"
$ws.Range("T4").WrapText = $true

$ws.Range("T5").Value = "No Mapping of such kind from Reg Processor`nID Repo- Might not be there in ID Repo as well"
$ws.Range("T5").WrapText = $true

$ws.Range("T6").Value = "ID Repo- need to know "

$ws.Range("T7").Value = "there shud be a label as Res_Service`nReg Client packet needs to be understood`nService from Reg proc needs to be developed"
$ws.Range("T7").WrapText = $true

# Existing cell S8 gets new text (and wraps now).
$ws.Range("S8").Value = "Reg proc`nArchival policy"
$ws.Range("S8").WrapText = $true

$ws.Range("T8").Value = "Under processing`nProcessed`n"
$ws.Range("T8").WrapText = $true

$ws.Range("T9").Value = "Under processing`nProcessed"
$ws.Range("T9").WrapText = $true

$ws.Range("T10").Value = "E-UIN Generation"

# --- Column width for the new column ------------------------------------
# (matches the saved worksheet width of ~32 characters, bestFit by the
# original author after typing the column's content)
$ws.Columns.Item(20).ColumnWidth = 31.09

# --- Selection: move the active cell in the frozen-pane view to T4 ------
$ws.Range("T4").Select()
